# Natmi following Dr Hou advice
# Recompute the Cx3cl1-Itgb3 LR-pairs table: the sender/target cluster set
# now includes "ECs", "FAPs", "sCs" (senders) crossed with "ECs", "FAPs",
# "M2", "sCs" (targets) for the Cx3cl1 -> Itgb3 ligand/receptor pair,
# replacing the old 3-cluster x 2 combinations with the new 3 x 4 grid
# and updated expression statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Cx3cl1", "Itgb3", "ECs",  2, 0.6666666666666666, 3.744854333333334, 11.234563,         0.2406267676902895, 0.2406267676902895, 3, 1,                  2.481489333333333,  7.444467999999999, 0.2345069082418988,  0.2345069082418987,  9.292816083053777,  83.635344747484,   0.0564286393312914,  0.05642863933129139),
    @("ECs",  "Cx3cl1", "Itgb3", "FAPs", 2, 0.6666666666666666, 3.744854333333334, 11.234563,         0.2406267676902895, 0.2406267676902895, 3, 1,                  7.245227,           21.735681,          0.6846919551326144,  0.6846919551326142,  27.13231972693367,  244.190877542403,  0.1647552120271057,  0.1647552120271056),
    @("ECs",  "Cx3cl1", "Itgb3", "M2",   2, 0.6666666666666666, 3.744854333333334, 11.234563,         0.2406267676902895, 0.2406267676902895, 2, 0.6666666666666666, 0.2001876666666667, 0.600563,           0.01891823194544989, 0.01891823194544989, 0.7496736509965556, 6.747062858969,    0.004552233003448784, 0.004552233003448783),
    @("ECs",  "Cx3cl1", "Itgb3", "sCs",  2, 0.6666666666666666, 3.744854333333334, 11.234563,         0.2406267676902895, 0.2406267676902895, 3, 1,                  0.6548283333333333, 1.964485,           0.06188290468003712, 0.06188290468003711, 2.452236721672778,  22.070130495055,   0.01489068332844362,  0.01489068332844361),
    @("FAPs", "Cx3cl1", "Itgb3", "ECs",  3, 1,                  10.317205,         30.951615,         0.662935182458301,  0.662935182458301,  3, 1,                  2.481489333333333,  7.444467999999999, 0.2345069082418988,  0.2345069082418987,  25.60203415731333,  230.4183074158199, 0.1554628800030752,  0.1554628800030752),
    @("FAPs", "Cx3cl1", "Itgb3", "FAPs", 3, 1,                  10.317205,         30.951615,         0.662935182458301,  0.662935182458301,  3, 1,                  7.245227,           21.735681,          0.6846919551326144,  0.6846919551326142,  74.750492230535,    672.7544300748149, 0.4539063862035705,  0.4539063862035704),
    @("FAPs", "Cx3cl1", "Itgb3", "M2",   3, 1,                  10.317205,         30.951615,         0.662935182458301,  0.662935182458301,  2, 0.6666666666666666, 0.2001876666666667, 0.600563,           0.01891823194544989, 0.01891823194544989, 2.065377195471666,  18.588394759245,   0.01254156154654528,  0.01254156154654528),
    @("FAPs", "Cx3cl1", "Itgb3", "sCs",  3, 1,                  10.317205,         30.951615,         0.662935182458301,  0.662935182458301,  3, 1,                  0.6548283333333333, 1.964485,           0.06188290468003712, 0.06188290468003711, 6.755998154808333,  60.803983393275,   0.04102435470511005,  0.04102435470511005),
    @("sCs",  "Cx3cl1", "Itgb3", "ECs",  3, 1,                  1.500857333333333, 4.502572,          0.09643804985140959, 0.09643804985140959, 3, 1,                 2.481489333333333,  7.444467999999999, 0.2345069082418988,  0.2345069082418987,  3.724361463521777,  33.51925317169599, 0.02261538890753217,  0.02261538890753216),
    @("sCs",  "Cx3cl1", "Itgb3", "FAPs", 3, 1,                  1.500857333333333, 4.502572,          0.09643804985140959, 0.09643804985140959, 3, 1,                 7.245227,           21.735681,          0.6846919551326144,  0.6846919551326142,  10.87405207461467,  97.86646867153199, 0.06603035690193816,  0.06603035690193815),
    @("sCs",  "Cx3cl1", "Itgb3", "M2",   3, 1,                  1.500857333333333, 4.502572,          0.09643804985140959, 0.09643804985140959, 2, 0.6666666666666666, 0.2001876666666667, 0.600563,           0.01891823194544989, 0.01891823194544989, 0.3004531275595555, 2.704078148036,    0.001824437395455826, 0.001824437395455826),
    @("sCs",  "Cx3cl1", "Itgb3", "sCs",  3, 1,                  1.500857333333333, 4.502572,          0.09643804985140959, 0.09643804985140959, 3, 1,                 0.6548283333333333, 1.964485,           0.06188290468003712, 0.06188290468003711, 0.9828039061577778, 8.845235155419999, 0.005967866646483448, 0.005967866646483446)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
